# Liêm Trinh Toạ Thủ và hội họp cát tinh
# Appends new Liêm Trinh Mệnh-palace rows (105-121) to Sheet1, columns A (title) and B (description).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 105-111: Liêm Trinh toạ thủ cung Mệnh tại từng cung (title in A, shared description in B)
$ws.Cells.Item(105, 1).Value = 'Liêm Trinh'
$ws.Cells.Item(105, 2).Value = 'Thân hình to lớn, xương thô, lông mày dầy.'
$ws.Cells.Item(106, 1).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Thìn'
$ws.Cells.Item(107, 1).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Tuất'
$ws.Cells.Item(106, 2).Value = 'Con người thẳng thắn, can đảm, dũng mãnh, nghiêm nghị, nóng tính.'
$ws.Cells.Item(109, 1).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Mùi'
$ws.Cells.Item(108, 1).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Sửu'
$ws.Cells.Item(110, 1).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Tý'
$ws.Cells.Item(111, 1).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Ngọ'
$ws.Cells.Item(107, 2).Value = 'Con người thẳng thắn, can đảm, dũng mãnh, nghiêm nghị, nóng tính.'
$ws.Cells.Item(108, 2).Value = 'Con người thẳng thắn, can đảm, dũng mãnh, nghiêm nghị, nóng tính.'
$ws.Cells.Item(109, 2).Value = 'Con người thẳng thắn, can đảm, dũng mãnh, nghiêm nghị, nóng tính.'
$ws.Cells.Item(110, 2).Value = 'Con người thẳng thắn, can đảm, dũng mãnh, nghiêm nghị, nóng tính.'
$ws.Cells.Item(111, 2).Value = 'Con người thẳng thắn, can đảm, dũng mãnh, nghiêm nghị, nóng tính.'

# Rows 112-121: Liêm Trinh toạ thủ cung Mệnh hội họp cát tinh (A and B identical per row)
$ws.Cells.Item(112, 1).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Tý gặp Quyền, Lộc, Khoa, Phủ, Tả, Hữu, Tướng, Xương, Khúc'
$ws.Cells.Item(112, 2).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Tý gặp Quyền, Lộc, Khoa, Phủ, Tả, Hữu, Tướng, Xương, Khúc'
$ws.Cells.Item(113, 1).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Ngọ gặp Quyền, Lộc, Khoa, Phủ, Tả, Hữu, Tướng, Xương, Khúc'
$ws.Cells.Item(113, 2).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Ngọ gặp Quyền, Lộc, Khoa, Phủ, Tả, Hữu, Tướng, Xương, Khúc'
$ws.Cells.Item(114, 1).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Dần gặp Quyền, Lộc, Khoa, Phủ, Tả, Hữu, Tướng, Xương, Khúc'
$ws.Cells.Item(114, 2).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Dần gặp Quyền, Lộc, Khoa, Phủ, Tả, Hữu, Tướng, Xương, Khúc'
$ws.Cells.Item(115, 1).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Thân gặp Quyền, Lộc, Khoa, Phủ, Tả, Hữu, Tướng, Xương, Khúc'
$ws.Cells.Item(115, 2).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Thân gặp Quyền, Lộc, Khoa, Phủ, Tả, Hữu, Tướng, Xương, Khúc'
$ws.Cells.Item(116, 1).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Thìn gặp Quyền, Lộc, Khoa, Phủ, Tả, Hữu, Tướng, Xương, Khúc'
$ws.Cells.Item(116, 2).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Thìn gặp Quyền, Lộc, Khoa, Phủ, Tả, Hữu, Tướng, Xương, Khúc'
$ws.Cells.Item(117, 1).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Tuất gặp Quyền, Lộc, Khoa, Phủ, Tả, Hữu, Tướng, Xương, Khúc'
$ws.Cells.Item(117, 2).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Tuất gặp Quyền, Lộc, Khoa, Phủ, Tả, Hữu, Tướng, Xương, Khúc'
$ws.Cells.Item(118, 1).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Sửu gặp Quyền, Lộc, Khoa, Phủ, Tướng, Xương, Khúc'
$ws.Cells.Item(118, 2).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Sửu gặp Quyền, Lộc, Khoa, Phủ, Tướng, Xương, Khúc'
$ws.Cells.Item(119, 1).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Mùi gặp Quyền, Lộc, Khoa, Phủ, Tướng, Xương, Khúc'
$ws.Cells.Item(119, 2).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Mùi gặp Quyền, Lộc, Khoa, Phủ, Tướng, Xương, Khúc'
$ws.Cells.Item(120, 1).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Sửu gặp Quyền, Lộc, Khoa, Phủ, Tả, Hữu, Tướng, Xương, Khúc'
$ws.Cells.Item(120, 2).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Sửu gặp Quyền, Lộc, Khoa, Phủ, Tả, Hữu, Tướng, Xương, Khúc'
$ws.Cells.Item(121, 1).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Mùi gặp Quyền, Lộc, Khoa, Phủ, Tả, Hữu, Tướng, Xương, Khúc'
$ws.Cells.Item(121, 2).Value = 'Liêm Trinh tọa thủ cung Mệnh ở Mùi gặp Quyền, Lộc, Khoa, Phủ, Tả, Hữu, Tướng, Xương, Khúc'

# Match the existing yellow highlight fill (style used by the rest of the sheet) for the new rows
$ws.Range("A105:B121").Interior.Color = 65535

# Leave the selection on the last edited cell, matching the saved view state
$ws.Range("B121").Select()
